# Applies the cryptos list update (GitHub Actions scheduled refresh).
# Updates Price (D) and Volume(1h) (E) columns for each coin row, plus a
# row swap for MultiversX / LidoDAOToken (rows 40-41 change order + values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.098.70'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').Value = '2.216.45'
$ws.Range('E3').Value = '  -1.19%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '''241.20'
$ws.Range('E5').Value = '  -2.07%  '
$ws.Range('E6').Value = '  -0.35%  '
$ws.Range('D7').Value = '''73.37'
$ws.Range('E7').Value = '  -1.50%  '
$ws.Range('E8').Value = '  +0.11%  '
$ws.Range('D9').Value = '''0.606'
$ws.Range('E9').Value = '  -1.44%  '
$ws.Range('D10').Value = '''42.94'
$ws.Range('E10').Value = '  +2.57%  '
$ws.Range('E11').Value = '  +1.15%  '
$ws.Range('D12').Value = '''7.10'
$ws.Range('E12').Value = '  +0.24%  '
$ws.Range('E13').Value = '  +0.32%  '
$ws.Range('D14').Value = '2.550.21'
$ws.Range('E14').Value = '  -1.15%  '
$ws.Range('D15').Value = '''14.19'
$ws.Range('E15').Value = '  -2.32%  '
$ws.Range('D16').Value = '''0.839'
$ws.Range('E16').Value = '  -1.34%  '
$ws.Range('D17').Value = '2.211.92'
$ws.Range('E17').Value = '  -0.86%  '
$ws.Range('D18').Value = '41.945.64'
$ws.Range('E18').Value = '  -0.39%  '
$ws.Range('E19').Value = '  +9.53%  '
$ws.Range('D20').Value = '''72.66'
$ws.Range('E20').Value = '  +0.96%  '
$ws.Range('D21').Value = '''6.13'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').Value = '''10.25'
$ws.Range('E22').Value = '  +17.22%  '
$ws.Range('D23').Value = '''229.68'
$ws.Range('E23').Value = '  -0.85%  '
$ws.Range('E24').Value = '  -7.27%  '
$ws.Range('D25').Value = '''11.66'
$ws.Range('E25').Value = '  +3.09%  '
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('D27').Value = '''3.59'
$ws.Range('E27').Value = '  -0.38%  '
$ws.Range('E28').Value = '  -2.02%  '
$ws.Range('E29').Value = '  +1.15%  '
$ws.Range('D30').Value = '''167.20'
$ws.Range('E30').Value = '  -1.15%  '
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('D32').Value = '''5.62'
$ws.Range('E32').Value = '  +8.89%  '
$ws.Range('D33').Value = '''0.0792'
$ws.Range('E33').Value = '  -2.89%  '
$ws.Range('E34').Value = '  +0.10%  '
$ws.Range('D35').Value = '''28.81'
$ws.Range('E35').Value = '  -4.69%  '
$ws.Range('E36').Value = '  -7.76%  '
$ws.Range('E37').Value = '  -4.95%  '
$ws.Range('D38').Value = '''0.0300'
$ws.Range('E38').Value = '  -2.82%  '
$ws.Range('D39').Value = '''13.15'
$ws.Range('E39').Value = '  -3.52%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').Value = '''2.12'
$ws.Range('E40').Value = '  -2.92%  '
$ws.Range('B41').Value = 'MultiversX'
$ws.Range('C41').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D41').Value = '''64.66'
$ws.Range('E41').Value = '  +4.00%  '
$ws.Range('E42').Value = '  -2.81%  '
$ws.Range('D43').Value = '''0.198'
$ws.Range('E43').Value = '  -2.17%  '
$ws.Range('D44').Value = '''8.73'
$ws.Range('E44').Value = '  +1.24%  '
$ws.Range('D45').Value = '''104.15'
$ws.Range('E45').Value = '  -2.38%  '
$ws.Range('E46').Value = '  -1.60%  '
$ws.Range('D47').Value = '''2.40'
$ws.Range('E47').Value = '  +6.21%  '
$ws.Range('E48').Value = '  -0.77%  '
$ws.Range('E49').Value = '  -0.28%  '
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('D51').Value = '2.421.76'
$ws.Range('E51').Value = '  -2.15%  '
